$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view state -------------------------------------------------
# The workbook was re-saved with the Excel window maximized (the classic
# "slightly negative origin / oversized" window rectangle Excel stores for
# a maximized window). Reproduce this as best effort via the COM window
# properties (values are in points; OOXML stores twips = points * 20).
$win = $excel.ActiveWindow
$win.WindowState = -4137   # xlMaximized
$win.Left   = -6
$win.Top    = -6
$win.Width  = 1932
$win.Height = 1056

# --- Row 6 (dps 12) now mirrors the layout used by rows 7/8 --------------
# Fill in the previously-empty "function"/"values" description for the
# dps-12 entry, matching the pattern used by the other rows, then copy the
# cell formatting (font/alignment) from row 7 so the row visually matches
# its neighbours.
$ws.Range("E6").Value = "0, 1"
$ws.Range("D6").Value = "isEmpty"

$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection -------------------------------------------------------------
# The saved cursor position moved to the Track-screen-related cell H25
# (reflecting that "Add Record" functionality moved there).
$ws.Range("H25").Select()
